$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new values look numeric; force Text format so Excel
# keeps them as literal strings instead of coercing to a Double.
$priceTextCells = [ordered]@{
    'D4' = '0.9998'
    'D5' = '244.43'
    'D6' = '0.6404'
    'D7' = '1.001'
    'D8' = '0.07489'
    'D9' = '0.2968'
    'D10' = '24.38'
    'D11' = '0.07654'
    'D13' = '5.032'
    'D14' = '0.6893'
    'D15' = '83.73'
    'D16' = '0.000009686'
    'D17' = '6.051'
    'D20' = '235.57'
    'D21' = '12.64'
    'D22' = '1.000'
    'D23' = '7.455'
    'D25' = '158.33'
    'D26' = '0.1412'
    'D27' = '8.518'
    'D28' = '17.91'
    'D29' = '0.06219'
    'D30' = '1.495'
    'D31' = '1.275'
    'D32' = '4.145'
    'D33' = '4.086'
    'D34' = '1.896'
    'D35' = '1.170'
    'D36' = '0.7269'
    'D37' = '2.608'
    'D38' = '2.831'
    'D39' = '0.01782'
    'D41' = '0.9220'
    'D42' = '6.130'
    'D43' = '1.001'
    'D45' = '102.17'
    'D47' = '0.00000000119'
    'D48' = '9.195'
    'D49' = '0.4053'
    'D51' = '1.644'
}
foreach ($addr in $priceTextCells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $priceTextCells[$addr]
    $rng.Style = "Normal"
}

# Remaining Price/Volume cells: plain text values (not numeric-ambiguous).
$plainCells = [ordered]@{
    'D2' = '29.678.53'
    'E2' = '  +1.54%  '
    'D3' = '1.854.37'
    'E3' = '  +1.16%  '
    'E5' = '  +0.82%  '
    'E6' = '  +3.07%  '
    'E7' = '  +0.04%  '
    'E8' = '  +1.61%  '
    'E9' = '  +1.94%  '
    'E10' = '  +4.90%  '
    'E11' = '  -0.22%  '
    'D12' = '1.864.88'
    'E12' = '  +1.91%  '
    'E13' = '  +1.16%  '
    'E14' = '  +3.01%  '
    'E15' = '  +1.24%  '
    'E16' = '  +8.06%  '
    'E17' = '  +3.09%  '
    'D18' = '29.703.99'
    'E18' = '  +1.71%  '
    'D19' = '2.107.91'
    'E19' = '  +1.96%  '
    'E20' = '  -0.25%  '
    'E21' = '  +1.13%  '
    'E22' = '  +0.03%  '
    'E23' = '  +1.36%  '
    'E24' = '  +0.08%  '
    'E25' = '  +0.12%  '
    'E26' = '  +0.60%  '
    'E27' = '  -0.28%  '
    'E28' = '  +1.48%  '
    'E29' = '  +8.07%  '
    'E30' = '  +0.53%  '
    'E31' = '  +5.24%  '
    'E32' = '  +1.36%  '
    'E33' = '  -0.48%  '
    'E34' = '  +1.60%  '
    'E35' = '  +2.38%  '
    'E36' = '  -0.21%  '
    'E37' = '  +0.24%  '
    'E38' = '  -1.05%  '
    'E39' = '  +1.47%  '
    'D40' = '1.201.03'
    'E40' = '  -1.87%  '
    'E41' = '  +1.70%  '
    'E42' = '  -2.32%  '
    'E43' = '  +0.02%  '
    'D44' = '2.016.70'
    'E44' = '  +2.25%  '
    'E45' = '  +0.59%  '
    'E46' = '  +1.38%  '
    'E47' = '  +1.25%  '
    'E48' = '  +0.74%  '
    'E49' = '  +0.70%  '
    'E50' = '  +1.03%  '
    'E51' = '  +3.10%  '
}
foreach ($addr in $plainCells.Keys) {
    $ws.Range($addr).Value = $plainCells[$addr]
}
